$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.008424997329712
$ws.Range("B1").Value = 2.118328809738159
$ws.Range("C1").Value = 6.428490161895752
$ws.Range("D1").Value = 1.635860562324524
$ws.Range("E1").Value = 1.366116642951965
